$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: extend the header row with two more columns (P1, Q1), copying the
# formatting of the existing header cell O1 (bold, centered, bordered style)
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: swap the I/K values and the M/O values, and append two new
# data columns P and Q (both valued 2) for every data row
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
